# Adds an "Instrument" worksheet to the workbook (common instrument fields
# factored out of CryoEMInstrument / XRayInstrument / SAXSInstrument), and
# adds a new "instruments" column to the Dataset sheet (referencing that
# new Instrument list), per the commit:
#   "Add instruments list to Dataset and change instrument_id to reference
#    Instrument objects"

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Instrument" worksheet right before "CryoEMInstrument" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Instrument"

$anchor = $wb.Worksheets.Item("CryoEMInstrument")
$newSheet.Move($anchor)

# Re-fetch by name: after Move() the old variable binding tracks the sheet
# that is now at the original tab position, not the moved sheet itself.
$ws = $wb.Worksheets.Item("Instrument")

# Match the outline summary placement used by every other sheet in the workbook.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

$ws.Range("A1").Value = "instrument_code"
$ws.Range("B1").Value = "manufacturer"
$ws.Range("C1").Value = "model"
$ws.Range("D1").Value = "installation_date"
$ws.Range("E1").Value = "current_status"
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "title"
$ws.Range("H1").Value = "description"

$ws.Range("E2:E1048576").Validation.Add(3, 1, 1, '"operational,maintenance,offline,commissioning"')
$ws.Range("E2:E1048576").Validation.IgnoreBlank = $true
$ws.Range("E2:E1048576").Validation.InCellDropdown = $true
$ws.Range("E2:E1048576").Validation.ShowInput = $false
$ws.Range("E2:E1048576").Validation.ShowError = $false

# --- 2. Add an "instruments" column to the Dataset sheet, before "studies" ---
$ds = $wb.Worksheets.Item("Dataset")
$ds.Columns.Item(2).Insert()
$ds.Range("B1").Value = "instruments"
